# Insert a new column before column D (this shifts D:K -> E:L and creates a fresh blank column D)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D1").EntireColumn.Insert()

# Copy number formats / styles from the (now shifted) column E into the new column D
# so the new column visually matches the rest of the year columns (date row uses the
# date format style, data rows use the numeric style).
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the FY2018 figures (period ending 2018-12-31 = serial 43465)
$rows = @(7,8,9,10,12,13,14,15,17,18,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,38,41,42,43,44,45,46,47,48,49,50,51,52,53,54,57,58,59,60,61,62,63,64,65,66,68,69,70,71,72,73,74,75,76,77,80,81,83,84,85,86,87,88,89,91,92,93,94,96,97,98,99,100,101,102)
$vals = @(43465,1531200,1364100,167100,"NA",0,0,1800,1480900,50300,-3600,86600,3700,43100,11800,0,31300,31100,0,"NA",0,0,3600,31100,0,31100,43465,7500,0,459300,"NA",8800,475600,18500,161900,89900,0,0,2900,0,748800,139500,4800,139500,283800,87600,52900,0,0,0,425800,0,0,0,0,174700,0,0,0,323000,0,43465,31100,39900,0,0,0,0,0,84800,-50700,0,0,-93200,0,0,0,0,10600,-100,2200)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $v = $vals[$i]
    $ws.Cells.Item($r, 4).Value2 = $v
}

$wb.Save()
